# Update confronto_modelli.xlsx with deeplab test results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cloud")

# Row 4 holds the Deeplabv3 results on the "cloud" sheet, previously "TODO"
# placeholders. Fill in the actual measured test results.

# Tempo inferenza [s] is a genuine number now.
$ws.Range("E4").Value = 89

# mIoU / WmIoU / min(IoU) / std(IoU): these cells are formatted as General but
# must keep storing their numeric-looking results as text (matching the other
# result cells in the table, e.g. "0.705", "0.684", ...). Assigning a numeric
# looking string directly via .Value would auto-convert the cell to a real
# number, so instead compute the text with TEXT() in a scratch cell and paste
# only the value, which preserves the text data type and the cell's style.
$scratch = $ws.Range("Z1")

$scratch.Formula = '=TEXT(0.746,"0.000")'
$scratch.Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4163) | Out-Null

$scratch.Formula = '=TEXT(0.728,"0.000")'
$scratch.Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4163) | Out-Null

$scratch.Formula = '=TEXT(0.457,"0.000")'
$scratch.Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4163) | Out-Null

$scratch.Formula = '=TEXT(0.166,"0.000")'
$scratch.Copy() | Out-Null
$ws.Range("J4").PasteSpecial(-4163) | Out-Null

$scratch.ClearContents() | Out-Null

# Move the active selection to J4 on the "cloud" sheet, matching the saved view.
$ws.Activate()
$ws.Range("J4").Select()
